$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sample's "Qty" header should read "qty" (lowercase) -- see commit
# "xlsx, csv, yaml: explode defaults to "qty" (#262) - Not "Qty"".
$ws.Range("B1").Value = "qty"

# Reflect the cell that was last edited/selected, matching the saved view state.
$ws.Range("B1").Select()
